$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3048.95356696664
$ws.Range("B3").Value = 3459.023402380415
$ws.Range("B4").Value = 3881.448748798804
$ws.Range("B5").Value = 4275.14330117802
$ws.Range("B6").Value = 4650.055980486404
$ws.Range("B7").Value = 4959.474869565301
$ws.Range("B8").Value = 5159.044660738143
$ws.Range("B9").Value = 5363.856978267961
$ws.Range("B10").Value = 5545.621323604955
$ws.Range("B11").Value = 5782.438858581801
$ws.Range("B12").Value = 5923.870091948377
$ws.Range("B13").Value = 6078.984671521705
$ws.Range("B14").Value = 6228.854537360922
$ws.Range("B15").Value = 6307.552041496362
$ws.Range("B16").Value = 6419.356968038338
$ws.Range("B17").Value = 6474.162623145135
$ws.Range("B18").Value = 6538.885201070315
$ws.Range("B19").Value = 6593.583980460076
$ws.Range("B20").Value = 6667.769515278936
$ws.Range("B21").Value = 6687.036915397022
$ws.Range("B22").Value = 6691.312506576869
$ws.Range("B23").Value = 6681.122391742525
$ws.Range("B24").Value = 6708.808511259751
$ws.Range("B25").Value = 6718.889982721571
$ws.Range("B26").Value = 6699.034090905601
$ws.Range("B27").Value = 6677.306450025442
$ws.Range("B28").Value = 6617.939411197557
$ws.Range("B29").Value = 6585.707603292922
$ws.Range("B30").Value = 6528.805143771848
$ws.Range("B31").Value = 6451.777864404964
$ws.Range("B32").Value = 6399.845128615047
$ws.Range("B33").Value = 6314.068930280756
$ws.Range("B34").Value = 6224.52253802735
$ws.Range("B35").Value = 6164.376395896112
$ws.Range("B36").Value = 6071.35755909473
$ws.Range("B37").Value = 5972.152035721722
$ws.Range("B38").Value = 5854.437045345408
$ws.Range("B39").Value = 5786.930234212736
$ws.Range("B40").Value = 5689.185453931108
$ws.Range("B41").Value = 5587.772366833796
$ws.Range("B42").Value = 5383.058315388869
$ws.Range("B43").Value = 5246.55293287376
$ws.Range("B44").Value = 5105.295375634791
$ws.Range("B45").Value = 4945.479565466244
$ws.Range("B46").Value = 4711.058470048256
$ws.Range("B47").Value = 4586.071295206078
$ws.Range("B48").Value = 4415.339258939557
$ws.Range("B49").Value = 4252.487202767622
$ws.Range("B50").Value = 4117.527047336341
$ws.Range("B51").Value = 3945.424246354755
$ws.Range("B52").Value = 3732.630817160296
$ws.Range("B53").Value = 3375.23342096731
$ws.Range("B54").Value = 3222.101485495768
$ws.Range("B55").Value = 3071.459351097113
$ws.Range("B56").Value = 2911.190209380306
$ws.Range("B57").Value = 2759.19461797026
$ws.Range("B58").Value = 2625.81891973139
$ws.Range("B59").Value = 2499.165860648438
$ws.Range("B60").Value = 2410.19285203388
$ws.Range("B61").Value = 2351.81319662543
$ws.Range("B62").Value = 2297.629950698137